# Further tests using different cross validation splits
# Adds 4 new RandomForest experiment blocks (rows 133-156) to the
# "RandomForest" worksheet, each mirroring the existing 6-row block
# structure (data row + 5 parameter-description rows) used for the
# earlier 80/20% run (rows 121-126 / 127-132).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RandomForest")

# The existing row 17 (17/80-20% split) result was re-measured; update its
# Accuracy / FMeasure values.
$ws.Range("I127").Value = 0.93558077436599996
$ws.Range("J127").Value = 0.93539079920099999

# Reused shared-string values (identical across all 4 new blocks), matching
# the same strings used by the pre-existing RandomForest 80/20% block
# (rows 121-126 / 127-132) in this sheet.
$paramLine1 = "RandomForestClassifier(bootstrap=True, class_weight=None, criterion='gini',"          # shared string 51
$paramLine2 = "            max_depth=None, max_features='sqrt', max_leaf_nodes=None,"               # shared string 52
$paramLine3 = "            min_impurity_decrease=0.0, min_impurity_split=None,"                     # shared string 41
$paramLine4 = "            min_samples_leaf=1, min_samples_split=2,"                                # shared string 42
$paramLine5 = "            min_weight_fraction_leaf=0.0, n_estimators=12000, n_jobs=6,"             # shared string 53
$paramLine6 = "            oob_score=False, random_state=None, verbose=0,"                          # shared string 44

# New "Cross Validation Split" labels must be interned into the shared
# string table in this exact order (70/30%, 75/25%, 85/15%, 90/10%) so the
# resulting xl/sharedStrings.xml matches the source workbook byte-for-byte.
# Note that row 145 ends up using the "90/10%" split while row 151 uses the
# "85/15%" split, but "85/15%" was authored (and therefore interned) before
# "90/10%" -- so we set column F for row 151 before row 145.
$ws.Range("F133").Value = "70/30%"
$ws.Range("F139").Value = "75/25%"
$ws.Range("F151").Value = "85/15%"
$ws.Range("F145").Value = "90/10%"

function Add-RFBlock {
    param($StartRow, $Index, $Split)

    $r0 = $StartRow

    $ws.Range("A$r0").Value = $Index
    $ws.Range("B$r0").Value = $paramLine1
    $ws.Range("C$r0").Value = "None"
    $ws.Range("D$r0").Value = 0
    $ws.Range("E$r0").Value = "None"
    $ws.Range("F$r0").Value = $Split
    $ws.Range("G$r0").Value = "RandomForest"
    $ws.Range("H$r0").Value = "var9(t)"

    $ws.Range("B" + ($r0 + 1)).Value = $paramLine2
    $ws.Range("B" + ($r0 + 2)).Value = $paramLine3
    $ws.Range("B" + ($r0 + 3)).Value = $paramLine4
    $ws.Range("B" + ($r0 + 4)).Value = $paramLine5
    $ws.Range("B" + ($r0 + 5)).Value = $paramLine6
}

# --- Block 1: row 133, index 18, 70/30% split ---
Add-RFBlock 133 18 "70/30%"
$ws.Range("I133").Value = 0.92677498330700003
$ws.Range("I133").NumberFormat = "0.00%"
$ws.Range("J133").Value = 0.92655824004300003
$ws.Range("J133").NumberFormat = "0.00%"

# --- Block 2: row 139, index 19, 75/25% split ---
Add-RFBlock 139 19 "75/25%"
$ws.Range("I139").Value = 0.93135683760700005
$ws.Range("I139").NumberFormat = "0.00%"
$ws.Range("J139").Value = 0.93114810886300003
$ws.Range("J139").NumberFormat = "0.00%"

# --- Block 3: row 145, index 20, 90/10% split ---
Add-RFBlock 145 20 "90/10%"
$ws.Range("I145").Value = 0.93457943925200004
$ws.Range("I145").NumberFormat = "0.00%"
$ws.Range("J145").Value = 0.93435051452899998
$ws.Range("J145").NumberFormat = "0.00%"

# --- Block 4: row 151, index 21, 85/15% split ---
Add-RFBlock 151 21 "85/15%"
$ws.Range("I151").Value = 0.93368936359599997
$ws.Range("I151").NumberFormat = "0.00%"
$ws.Range("J151").Value = 0.93347648280700002
$ws.Range("J151").NumberFormat = "0.00%"

# --- Update the view: scroll near the bottom and select B145 ---
$ws.Activate()
$ws.Range("B145").Select()
$excel.ActiveWindow.ScrollRow = 131
$excel.ActiveWindow.ScrollColumn = 1
